$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# Row 13: Rubeshen
$row13 = $tbl.ListRows.Add()
$row13.Range.Cells.Item(1, 1).Value = "Rubeshen"
$row13.Range.Cells.Item(1, 2).Value = "23:35 15/05/1992 +08:00"
$row13.Range.Cells.Item(1, 3).Value = "Klang"
$row13.Range.Cells.Item(1, 4).Value = 101.4456
$row13.Range.Cells.Item(1, 5).Value = 3.0449000000000002
$row13.Range.Cells.Item(1, 6).Value = "Male"
$row13.Range.Cells.Item(1, 6).NumberFormat = "@"

# Row 14: Haravin
$row14 = $tbl.ListRows.Add()
$row14.Range.Cells.Item(1, 1).Value = "Haravin"
$row14.Range.Cells.Item(1, 2).Value = "04:10 25/03/1991 +08:00"
$row14.Range.Cells.Item(1, 3).Value = "Kuala Lumpur"
$row14.Range.Cells.Item(1, 4).Value = 101.68729999999999
$row14.Range.Cells.Item(1, 5).Value = 3.1602700000000001
$row14.Range.Cells.Item(1, 6).Value = "Male"

# View state: zoom to 130% and move the selection to C11
$excel.ActiveWindow.Zoom = 130
$ws.Range("C11").Select() | Out-Null
